$d = $word.ActiveDocument

function Set-QuotedBoldUnderline($fullQuote, $innerText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($fullQuote, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find: $fullQuote"
    }
    $prefixLen = "【問】『".Length
    $suffixLen = "』".Length
    $s = $rng.Start
    $e = $rng.End
    $inner = $d.Range($s + $prefixLen, $e - $suffixLen)
    if ($inner.Text -ne $innerText) {
        throw "Mismatch: expected [$innerText] got [$($inner.Text)]"
    }
    $inner.Font.Bold = $true
    $inner.Font.Underline = 1
}

Set-QuotedBoldUnderline "【問】『開善』" "開善"

Set-QuotedBoldUnderline "【問】『能於黑色，通達一切，非於一切。非通達一切 ，是通達一切，非非非是。一切法邪，一切法正』" "能於黑色，通達一切，非於一切。非通達一切 ，是通達一切，非非非是。一切法邪，一切法正"

Set-QuotedBoldUnderline "【問】『不斷生身菩薩之近疑』" "不斷生身菩薩之近疑"

Write-Host "done"
